$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 407, shifting rows 407:463 down to 408:464.
$ws.Rows.Item(407).Insert()

# Populate the newly inserted row 407 with the latest weekly reading.
$ws.Cells.Item(407, 1).Value = 9
$ws.Cells.Item(407, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(407, 3).Value = "Metropolitana"
$ws.Cells.Item(407, 4).Value = 45127
$ws.Cells.Item(407, 5).Value = 13
$ws.Cells.Item(407, 6).Value = 300000001
$ws.Cells.Item(407, 7).Value = "Rabanito"
$ws.Cells.Item(407, 8).Value = "Sin especificar"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 7000
$ws.Cells.Item(407, 11).Value = 3000
$ws.Cells.Item(407, 12).Value = 4000
$ws.Cells.Item(407, 13).Value = 3500
$ws.Cells.Item(407, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(407, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(407, 16).Value = 35
$ws.Cells.Item(407, 17).Value = 100
$ws.Cells.Item(407, 18).Value = "Hortaliza"
